# Add two new data rows under the existing "Usuario / Contraseña / Tipo"
# header row. The Usuario/Contraseña columns hold numeric-looking login
# codes ("123", "2") that must be stored as TEXT, not numbers, so they are
# entered with a leading apostrophe (the normal Excel "force text" input),
# matching how a user would type them in the grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'123"
$ws.Range("B2").Value = "'123"
$ws.Range("C2").Value = "Cliente"

$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "'2"
$ws.Range("C3").Value = "Cliente"

# The apostrophe-prefix entry flags the cells with the "quote prefix" text
# style; reset back to Normal so the new rows keep the default (unstyled)
# look, same as the header row keeps its own existing style untouched.
$ws.Range("A2:B3").Style = "Normal"
